$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sample_id values in column C: insert "." between the colony number
# and the trailing "2" (second fragment marker), e.g. "H32C" -> "H3.2C".
# This corresponds to the newly added fvfm metadata notation in the commit.
$ws.Range("C2").Value = "H3.2C"
$ws.Range("C3").Value = "H2.2C"
$ws.Range("C4").Value = "H1.2C"
$ws.Range("C5").Value = "H10.2C"
$ws.Range("C6").Value = "H8.2C"
$ws.Range("C7").Value = "A6.2H"
$ws.Range("C8").Value = "A9.2H"
$ws.Range("C9").Value = "A10.2H"
$ws.Range("C10").Value = "A3.2H"
$ws.Range("C11").Value = "A7.2H"
$ws.Range("C12").Value = "H5.2L"
$ws.Range("C13").Value = "H2.2L"
$ws.Range("C14").Value = "H3.2L"
$ws.Range("C15").Value = "H9.2L"
$ws.Range("C16").Value = "H8.2L"
$ws.Range("C17").Value = "A2.2C"
$ws.Range("C18").Value = "A9.2C"
$ws.Range("C19").Value = "A10.2C"
$ws.Range("C20").Value = "A7.2C"
$ws.Range("C21").Value = "A1.2C"
$ws.Range("C22").Value = "H6.2L"
$ws.Range("C23").Value = "H4.2L"
$ws.Range("C24").Value = "H10.2L"
$ws.Range("C25").Value = "H7.2L"
$ws.Range("C26").Value = "H1.2L"
$ws.Range("C27").Value = "A6.2C"
$ws.Range("C28").Value = "A4.2C"
$ws.Range("C29").Value = "A8.2C"
$ws.Range("C30").Value = "A5.2C"
$ws.Range("C31").Value = "A3.2C"
$ws.Range("C32").Value = "H7.2H"
$ws.Range("C33").Value = "H8.2H"
$ws.Range("C34").Value = "H10.2H"
$ws.Range("C35").Value = "H1.2H"
$ws.Range("C36").Value = "H4.2H"
$ws.Range("C37").Value = "A9.2L"
$ws.Range("C38").Value = "A6.2L"
$ws.Range("C39").Value = "A10.2L"
$ws.Range("C40").Value = "A4.2L"
$ws.Range("C41").Value = "A3.2L"
$ws.Range("C42").Value = "H6.2H"
$ws.Range("C43").Value = "H2.2H"
$ws.Range("C44").Value = "H5.2H"
$ws.Range("C45").Value = "H3.2H"
$ws.Range("C46").Value = "H9.2H"
$ws.Range("C47").Value = "A8.2L"
$ws.Range("C48").Value = "A1.2L"
$ws.Range("C49").Value = "A2.2L"
$ws.Range("C50").Value = "A7.2L"
$ws.Range("C51").Value = "A5.2L"
$ws.Range("C52").Value = "H9.2C"
$ws.Range("C53").Value = "H6.2C"
$ws.Range("C54").Value = "H5.2C"
$ws.Range("C55").Value = "H7.2C"
$ws.Range("C56").Value = "H4.2C"
$ws.Range("C57").Value = "A5.2H"
$ws.Range("C58").Value = "A4.2H"
$ws.Range("C59").Value = "A8.2H"
$ws.Range("C60").Value = "A2.2H"
$ws.Range("C61").Value = "A1.2H"

# Update the active selection to match the saved view state (C62).
$ws.Range("C62").Select()
